# fix: ajustar para preencher dados a partir da linha 2 (A2) mantendo cabeçalho na linha 1
#
# The template previously shipped with two sample/example rows of data
# (row 2 "Contato 1" and row 3 "Contato 2") including mailto: hyperlinks on
# column B. This cleared those sample rows so the template starts empty
# right after the header row, while keeping the header row (row 1) and the
# formatting/styles of the first two data rows intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the mailto hyperlinks that were attached to the sample e-mail
# addresses in B2 and B3.
$ws.Range("B2").Hyperlinks.Delete()
$ws.Range("B3").Hyperlinks.Delete()

# Clear the sample data in rows 2 and 3 (columns A through P) but keep the
# cell formatting/styles so the sheet still looks consistent for new data.
$ws.Range("A2:P3").ClearContents()

# Move the active selection/view back to the top-left of the sheet, now
# focused on the first empty data cell (B7 per the original author's last
# selection) instead of the old N1 selection with C1 scrolled into view.
$ws.Range("B7").Select()
